# Regenerate the "K" column (column G) values in the save_data sheet.
# The source data behind this sheet was regenerated upstream (K replaces
# the old Strike# derived values), so we overwrite column G with the
# newly computed strikeout counts for each game row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> new K (column G) value.
$newValues = @{
    2  = 1
    3  = 0
    4  = 2
    5  = 0
    6  = 0
    7  = 0
    8  = 2
    9  = 1
    10 = 0
    11 = 1
    12 = 1
    13 = 0
    14 = 1
    15 = 2
    16 = 1
    17 = 0
    18 = 0
    19 = 0
    20 = 0
    21 = 2
    22 = 2
    23 = 1
    24 = 0
    25 = 1
    26 = 1
    27 = 2
    28 = 0
    29 = 1
    31 = 0
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
